# Add new enforcement-action rows' details (columns D:M) for rows 55-58.
# New shared strings must be created in the exact order they first appear so
# that the resulting sharedStrings.xml table indices line up with the target
# workbook. That order (as created by the original author) is:
#   1) ICOS   2) ICOBox   3) Bitqyck Inc.   4) Unregistered Offering and Exchange   5) HLTH
# which means H56 ("Bitqyck Inc.") must be written before E56
# ("Unregistered Offering and Exchange").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 55: SEC v. ICOBox, et al. ---
$ws.Range("D55").Value = "Ongoing"
$ws.Range("E55").Value = "Unregistered Offering"
$ws.Range("F55").Value = "Civil"
$ws.Range("G55").Value = "ICOS"
$ws.Range("H55").Value = "ICOBox"
$ws.Range("I55").Value = "Ethereum"
$ws.Range("J55").Value = 14600000
$ws.Range("K55").Value = 1
$ws.Range("L55").Value = 1
$ws.Range("M55").Value = "Los Angeles"

# --- Row 56: SEC v. Bitqyck, Inc., et al. ---
$ws.Range("D56").Value = "Settlement"
$ws.Range("H56").Value = "Bitqyck Inc."
$ws.Range("E56").Value = "Unregistered Offering and Exchange"
$ws.Range("F56").Value = "Civil"
$ws.Range("G56").Value = "N/A"
$ws.Range("I56").Value = "Ethereum"
$ws.Range("J56").Value = 13000000
$ws.Range("K56").Value = 1
$ws.Range("L56").Value = 1
$ws.Range("M56").Value = "Fort Worth"

# --- Row 57: ICO Rating cease-and-desist ---
$ws.Range("D57").Value = "Settlement"
$ws.Range("E57").Value = "Anti-touting"
$ws.Range("F57").Value = "Civil"
$ws.Range("G57").Value = "N/A"
$ws.Range("H57").Value = "ICO Rating"
$ws.Range("I57").Value = "N/A"
$ws.Range("J57").Value = 268998
$ws.Range("K57").Value = 1
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = "Washington, D.C."

# --- Row 58: SimplyVital Health, Inc. ---
$ws.Range("D58").Value = "Settlement"
$ws.Range("E58").Value = "Unregistered Offering"
$ws.Range("F58").Value = "Civil"
$ws.Range("G58").Value = "HLTH"
$ws.Range("H58").Value = "SimplyVital Health, Inc."
$ws.Range("I58").Value = "Ethereum"
$ws.Range("J58").Value = 6300000
$ws.Range("K58").Value = 1
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = "Boston"

# Widen column E slightly (author manually resized it) and move the
# selection to F58, matching the cursor position left by the edit.
$ws.Columns.Item(5).ColumnWidth = 30.833333333333332
$ws.Range("F58").Select()
